$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header cells (AC1:AE1) for the team record columns ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the look of the rest of the header row (bold/centered/bordered)
# by copying the formatting from the neighboring header cell (AB1).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the team's season record (66 wins - 49 losses - 0 ties) ---
# for every player row in the sheet (rows 2 through 35).
$lastRow = 35
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 66   # column AC - Wins
    $ws.Cells.Item($r, 30).Value = 49   # column AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # column AE - Ties
}
